$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45186 -> 2023-09-17) that needs
# to be bumped to 45188 (-> 2023-09-19) for every data row (rows 2 through 400).
$lastRow = 400

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
